# Added -Z end cap connector to 1U
# - Bump the PCB version note
# - Insert a new BOM line for J1 (the -Z end cap connector)
# - Record the new revision history entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PCB version note
$ws.Range("A3").Value = "PCB version: 1.1"

# The "Insert" range method isn't available in this host, so emulate
# "insert a blank row at 12" by shifting the trailing rows (13, 15, 16, 17)
# down one row each, working from the bottom up so nothing gets clobbered.
# Destinations are cleared first since Copy() leaves a cell alone when the
# matching source cell is blank.
$ws.Range("A18:J18").ClearContents()
$ws.Range("A17:J17").Copy($ws.Range("A18:J18"))
$ws.Rows(18).RowHeight = 13.5

$ws.Range("A17:J17").ClearContents()
$ws.Range("A16:J16").Copy($ws.Range("A17:J17"))
$ws.Rows(17).RowHeight = 13.5

$ws.Range("A16:J16").ClearContents()
$ws.Range("A15:J15").Copy($ws.Range("A16:J16"))
$ws.Rows(16).RowHeight = 13.5

$ws.Range("A14:J14").ClearContents()
$ws.Range("A13:J13").Copy($ws.Range("A14:J14"))
$ws.Rows(14).RowHeight = 13.5

# Blank out the rows that moved away from their old positions (13 and 15
# no longer hold anything -- and since nothing refers to them anymore,
# drop their custom height too so no empty row lingers in the XML).
$ws.Range("A13:J13").ClearContents()
$ws.Rows(13).AutoFit()
$ws.Range("A15:J15").ClearContents()
$ws.Rows(15).AutoFit()

# New BOM row for connector J1, styled like the other part rows (row 11)
$ws.Range("A12:J12").ClearContents()
$ws.Range("A11:J11").Copy($ws.Range("A12:J12"))
$ws.Rows(12).RowHeight = 13.5

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "J1"
$ws.Range("C12").Value = "P"
$ws.Range("D12").Value = "TE"
$ws.Range("E12").Value = "5-104196-5"
$ws.Range("F12").Value = "20 Position Receptacle Connector  Through Hole, Right Angle"
$ws.Range("G12").Value = "Digi-Key"
$ws.Range("H12").Value = "A115240-ND"
$ws.Range("I12").Value = "N"
$ws.Range("J12").ClearContents()

# New revision-history row
$ws.Range("A18").Value = "1.1r0"
$ws.Range("B18").Value = "Added -Z end cap connector"

# Restore the active selection cell noted in the saved view state
$ws.Range("J12").Select()
